$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "ingrid.matos@mrv.com.br"
$ws.Range("B4").Value = "PowerBI"
$ws.Range("C4").Value = "Painel Power BI"
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = "muito importante também"
